# PRJ0018886_Hierarchy Viewer+ Time recordManager(PArtial changes)
#
# Updates the "Project Clear" / "Project Bend" engagement names across the
# lookup sheets, adds a 4th (helper) column to Project_Title with the full
# "Project Bend" name, bolds the Project_Title header row, and refreshes
# the active sheet/selection bookkeeping.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Project_Title sheet: header row bold, new values + helper column D
# ---------------------------------------------------------------------
$wsProjectTitle = $wb.Worksheets.Item("Project_Title")

$wsProjectTitle.Range("A1:B1").Font.Bold = $true

$wsProjectTitle.Range("D2").Value = "Project Bend-Bernhard Capital Partners Management-FVA-109081"
$wsProjectTitle.Range("A2").Value = "Project Clear-LucidHealth-FVA-105379"
$wsProjectTitle.Range("B2").Value = "Fieldwork"

$wsProjectTitle.Columns.Item(1).ColumnWidth = 55.6640625
$wsProjectTitle.Columns.Item(2).ColumnWidth = 15
$wsProjectTitle.Columns.Item(4).ColumnWidth = 55.6640625

$wsProjectTitle.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# RateSheetManagement sheet: engagement name refresh
# ---------------------------------------------------------------------
$wsRateSheet = $wb.Worksheets.Item("RateSheetManagement")
$wsRateSheet.Range("A2").Value = "Project Clear-LucidHealth-FVA-105379"
$wsRateSheet.Range("C19").Select() | Out-Null

# ---------------------------------------------------------------------
# WeeklyEntryMatrix sheet: engagement name refresh
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("WeeklyEntryMatrix")
$wsWeekly.Range("A2").Value = "Project Clear-LucidHealth-FVA-105379"
$wsWeekly.Columns.Item(1).ColumnWidth = 24.77734375

# ---------------------------------------------------------------------
# Activate Project_Title and select D9 (new last-saved view state)
# ---------------------------------------------------------------------
$wsProjectTitle.Activate()
$wsProjectTitle.Range("D9").Select() | Out-Null
